$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Stephens-MacCall filtered" sample-size columns (C = the
# first/combined filter, D = the drift-level filter) with new values,
# ahead of removing the drift-level SM filtering.
# Write column D first, then column C, so new shared-string entries are
# appended in the same order the original authoring tool produced them.
$ws.Range("D2").Value = "706 (68%)"
$ws.Range("D3").Value = "1813 (91%)"
$ws.Range("D4").Value = "806 (62%)"
$ws.Range("D5").Value = "798 (57%)"
$ws.Range("D6").Value = "1449 (81%)"
$ws.Range("D7").Value = "1627 (85%)"
$ws.Range("C2").Value = "3038 (30%)"
$ws.Range("C3").Value = "7490 (60%)"
$ws.Range("C4").Value = "2740 (31%)"
$ws.Range("C5").Value = "1331 (22%)"
$ws.Range("C6").Value = "5088 (45%)"
$ws.Range("C7").Value = "5040 (45%)"

$ws.Columns.Item(2).ColumnWidth = 15.6
$ws.Columns.Item(3).ColumnWidth = 22.0

# Mirror the saved cursor position recorded in the sheet's selection.
$ws.Range("C5").Select()
